# The commit swaps the two embedded themes of the deck: the theme bound to
# the (only) slide master - "Integral" - becomes the stock "Office Theme"
# palette, while the palette previously shipped as the secondary theme
# ("Office Theme", used by the notes master) becomes "Integral".
#
# Font scheme / format scheme (fills, lines, effects) are identical between
# the two themes, so the only real content change is the 12-slot theme
# color scheme (and, cosmetically, the <a:theme>/<a:clrScheme> name
# attributes, which PowerPoint does not expose as a writable property).
#
# We apply the reachable half of that swap through the documented
# PowerPoint object model: Master.Theme.ThemeColorScheme exposes all 12
# theme colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) as
# read/write RGB values, which is exactly the mechanism PowerPoint itself
# uses when a user recolors a theme.

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$themeColors = $master.Theme.ThemeColorScheme

# Target palette: the stock "Office Theme" color scheme.
# Order matches the 12 DrawingML clrScheme slots:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeRGB = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $officeThemeRGB.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeRGB[$i - 1]
}
